# Auto-generated edit script: refresh cached market-price figures
# (currentAveragePrice / LevePrice / LeveProfit columns) per sheet,
# matching the scheduled-runner commit that re-pulled Universalis data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 634438.4
$ws.Range("J17").Value = 654734.75
$ws.Range("L17").Value = 1964204.25
$ws.Range("N17").Value = -1964540.25
$ws.Range("H69").Value = 7716.25
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = ""
$ws.Range("H72").Value = 7716.25
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = ""
$ws.Range("H76").Value = 3975
$ws.Range("I76").Value = 3701.5
$ws.Range("K76").Value = 3701.5
$ws.Range("M76").Value = -3386.5
$ws.Range("H79").Value = 3975
$ws.Range("I79").Value = 3701.5
$ws.Range("K79").Value = 3701.5
$ws.Range("M79").Value = -2609.5
$ws.Range("H101").Value = 1068.8182
$ws.Range("I101").Value = 977.2857
$ws.Range("J101").Value = 1229
$ws.Range("K101").Value = 2931.8571
$ws.Range("L101").Value = 3687
$ws.Range("M101").Value = -1309.8571
$ws.Range("N101").Value = -6931
$ws.Range("H103").Value = 580.46155
$ws.Range("J103").Value = 792.1667
$ws.Range("L103").Value = 2376.5001
$ws.Range("N103").Value = -3548.5001
$ws.Range("H123").Value = 125000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 125000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 125000
$ws.Range("M123").Value = ""
$ws.Range("N123").Value = -134800
$ws.Range("H127").Value = 891.3889
$ws.Range("I127").Value = 736.4
$ws.Range("K127").Value = 2209.2
$ws.Range("M127").Value = 2750.8
$ws.Range("H129").Value = 1780.5385
$ws.Range("J129").Value = 1916.6666
$ws.Range("L129").Value = 5749.9998
$ws.Range("N129").Value = -15749.9998
$ws.Range("H132").Value = 1831.697
$ws.Range("I132").Value = 1885.1613
$ws.Range("K132").Value = 5655.4839
$ws.Range("M132").Value = -3125.4839
$ws.Range("H138").Value = 3037.54
$ws.Range("I138").Value = 2148.2307
$ws.Range("K138").Value = 6444.6921
$ws.Range("M138").Value = -1304.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3041.182
$ws.Range("I45").Value = 2494.8147
$ws.Range("K45").Value = 2494.8147
$ws.Range("M45").Value = -2117.8147
$ws.Range("H102").Value = 4625
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 5833.3335
$ws.Range("K102").Value = 1000
$ws.Range("L102").Value = 5833.3335
$ws.Range("M102").Value = 622
$ws.Range("N102").Value = -9077.333500000001
$ws.Range("H110").Value = 3499.6667
$ws.Range("I110").Value = 3499.5
$ws.Range("K110").Value = 3499.5
$ws.Range("M110").Value = -1454.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 245.6
$ws.Range("J22").Value = 245.75
$ws.Range("L22").Value = 245.75
$ws.Range("N22").Value = -591.75
$ws.Range("H105").Value = 1544.5405
$ws.Range("I105").Value = 1454.3823
$ws.Range("K105").Value = 1454.3823
$ws.Range("M105").Value = 292.6177
$ws.Range("H107").Value = 2395.1333
$ws.Range("I107").Value = 1893.625
$ws.Range("K107").Value = 1893.625
$ws.Range("M107").Value = 26.375
$ws.Range("H141").Value = 188248.5
$ws.Range("J141").Value = 188248.5
$ws.Range("L141").Value = 188248.5
$ws.Range("N141").Value = -198608.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8624.929
$ws.Range("I31").Value = 1400
$ws.Range("K31").Value = 1400
$ws.Range("M31").Value = -1105
$ws.Range("H34").Value = 8624.929
$ws.Range("I34").Value = 1400
$ws.Range("K34").Value = 1400
$ws.Range("M34").Value = -1198
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("H68").Value = 47998.453
$ws.Range("J68").Value = 51664.777
$ws.Range("L68").Value = 51664.777
$ws.Range("N68").Value = -53162.777
$ws.Range("H71").Value = 47998.453
$ws.Range("J71").Value = 51664.777
$ws.Range("L71").Value = 154994.331
$ws.Range("N71").Value = -162482.331
$ws.Range("H99").Value = 3050.5833
$ws.Range("I99").Value = 3024.6
$ws.Range("J99").Value = 3180.5
$ws.Range("K99").Value = 3024.6
$ws.Range("L99").Value = 3180.5
$ws.Range("M99").Value = -1526.6
$ws.Range("N99").Value = -6176.5
$ws.Range("H126").Value = 3050.5833
$ws.Range("I126").Value = 3024.6
$ws.Range("J126").Value = 3180.5
$ws.Range("K126").Value = 9073.799999999999
$ws.Range("L126").Value = 9541.5
$ws.Range("M126").Value = -6603.799999999999
$ws.Range("N126").Value = -14481.5
$ws.Range("H132").Value = 2895.7568
$ws.Range("I132").Value = 2683.647
$ws.Range("K132").Value = 8050.941
$ws.Range("M132").Value = -5520.941
$ws.Range("H133").Value = 90000
$ws.Range("J133").Value = 90000
$ws.Range("L133").Value = 90000
$ws.Range("N133").Value = -95060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2.2
$ws.Range("J12").Value = 2.2
$ws.Range("L12").Value = 6.600000000000001
$ws.Range("N12").Value = -352.6
$ws.Range("H113").Value = 1341.2858
$ws.Range("I113").Value = 507.5
$ws.Range("J113").Value = 1674.8
$ws.Range("K113").Value = 1522.5
$ws.Range("L113").Value = 5024.4
$ws.Range("N113").Value = -9364.4
$ws.Range("M113").Value = 647.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 30999.75
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("H120").Value = 66945
$ws.Range("J120").Value = 66945
$ws.Range("L120").Value = 66945
$ws.Range("N120").Value = -76621
$ws.Range("H123").Value = 39000
$ws.Range("J123").Value = 39000
$ws.Range("L123").Value = 39000
$ws.Range("N123").Value = -43900
$ws.Range("H133").Value = 112999.5
$ws.Range("J133").Value = 112999.5
$ws.Range("L133").Value = 112999.5
$ws.Range("N133").Value = -123119.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6223.7144
$ws.Range("I7").Value = 6224.273
$ws.Range("J7").Value = 6221.6665
$ws.Range("K7").Value = 6224.273
$ws.Range("L7").Value = 6221.6665
$ws.Range("M7").Value = -6112.273
$ws.Range("N7").Value = -6445.6665
$ws.Range("H40").Value = 19611802
$ws.Range("I40").Value = 23812152
$ws.Range("K40").Value = 23812152
$ws.Range("M40").Value = -23812016
$ws.Range("H126").Value = 6223.7144
$ws.Range("I126").Value = 6224.273
$ws.Range("J126").Value = 6221.6665
$ws.Range("K126").Value = 18672.819
$ws.Range("L126").Value = 18664.9995
$ws.Range("M126").Value = -16202.819
$ws.Range("N126").Value = -23604.9995
$ws.Range("H136").Value = 3514
$ws.Range("I136").Value = 2199.6667
$ws.Range("J136").Value = 4499.75
$ws.Range("K136").Value = 6599.000100000001
$ws.Range("L136").Value = 13499.25
$ws.Range("M136").Value = -4049.000100000001
$ws.Range("N136").Value = -18599.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3573.8462
$ws.Range("I81").Value = 2725.2942
$ws.Range("J81").Value = 5176.6665
$ws.Range("K81").Value = 5450.5884
$ws.Range("L81").Value = 10353.333
$ws.Range("M81").Value = -4389.5884
$ws.Range("N81").Value = -12475.333
$ws.Range("H84").Value = 3573.8462
$ws.Range("I84").Value = 2725.2942
$ws.Range("J84").Value = 5176.6665
$ws.Range("K84").Value = 27252.942
$ws.Range("L84").Value = 51766.665
$ws.Range("M84").Value = -21948.942
$ws.Range("N84").Value = -62374.665
$ws.Range("H132").Value = 2774.3242
$ws.Range("I132").Value = 2401.5
$ws.Range("J132").Value = 6999.6665
$ws.Range("K132").Value = 7204.5
$ws.Range("L132").Value = 20998.9995
$ws.Range("M132").Value = -4674.5
$ws.Range("N132").Value = -26058.9995
$ws.Range("H136").Value = 25791.428
$ws.Range("J136").Value = 94431.73
$ws.Range("L136").Value = 283295.19
$ws.Range("N136").Value = -288395.19
